$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2873.6924
$ws.Range("I32").Value = 3112.5
$ws.Range("K32").Value = 3112.5
$ws.Range("M32").Value = -2786.5
$ws.Range("H43").Value = 8171.1113
$ws.Range("I43").Value = 6398.5
$ws.Range("K43").Value = 6398.5
$ws.Range("M43").Value = -6329.5
$ws.Range("H62").Value = 71592
$ws.Range("I62").Value = 95035.27
$ws.Range("K62").Value = 95035.27
$ws.Range("M62").Value = -94411.27
$ws.Range("H65").Value = 71592
$ws.Range("I65").Value = 95035.27
$ws.Range("K65").Value = 475176.35
$ws.Range("M65").Value = -472056.35
$ws.Range("H100").Value = 3101.1538
$ws.Range("I100").Value = 2900.375
$ws.Range("K100").Value = 2900.375
$ws.Range("M100").Value = -2359.375
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H111").Value = 2098.3333
$ws.Range("I111").Value = 2098.3333
$ws.Range("K111").Value = 6294.999899999999
$ws.Range("M111").Value = -3227.999899999999
$ws.Range("H138").Value = 2431.1516
$ws.Range("J138").Value = 3523.5334
$ws.Range("L138").Value = 10570.6002
$ws.Range("N138").Value = -20850.6002

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 162.9
$ws.Range("I5").Value = 139.58333
$ws.Range("K5").Value = 139.58333
$ws.Range("M5").Value = -27.58332999999999
$ws.Range("H32").Value = 37299.137
$ws.Range("I32").Value = 20771.725
$ws.Range("J32").Value = 157713.14
$ws.Range("K32").Value = 20771.725
$ws.Range("L32").Value = 157713.14
$ws.Range("M32").Value = -20484.725
$ws.Range("N32").Value = -158287.14
$ws.Range("H55").Value = 19916.666
$ws.Range("H63").Value = 1970.6428
$ws.Range("J63").Value = 2399
$ws.Range("L63").Value = 2399
$ws.Range("N63").Value = -3771
$ws.Range("H66").Value = 1970.6428
$ws.Range("J66").Value = 2399
$ws.Range("L66").Value = 11995
$ws.Range("N66").Value = -18859
$ws.Range("H88").Value = 13334113
$ws.Range("J88").Value = 22222894
$ws.Range("L88").Value = 22222894
$ws.Range("N88").Value = -22223706
$ws.Range("H91").Value = 13334113
$ws.Range("J91").Value = 22222894
$ws.Range("L91").Value = 22222894
$ws.Range("N91").Value = -22225702
$ws.Range("H110").Value = 3343.8572
$ws.Range("I110").Value = 3102
$ws.Range("J110").Value = 3666.3333
$ws.Range("K110").Value = 3102
$ws.Range("L110").Value = 3666.3333
$ws.Range("M110").Value = -1057
$ws.Range("N110").Value = -7756.3333
$ws.Range("H122").Value = 1821.6296
$ws.Range("I122").Value = 1766.84
$ws.Range("K122").Value = 5300.52
$ws.Range("M122").Value = -2850.52

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 162.9
$ws.Range("I4").Value = 139.58333
$ws.Range("K4").Value = 139.58333
$ws.Range("M4").Value = -24.58332999999999
$ws.Range("H19").Value = 900000000
$ws.Range("J19").Value = 900000000
$ws.Range("L19").Value = 900000000
$ws.Range("N19").Value = -900000346
$ws.Range("H82").Value = 22248.555
$ws.Range("J82").Value = 24997.143
$ws.Range("L82").Value = 24997.143
$ws.Range("N82").Value = -25763.143
$ws.Range("H85").Value = 22248.555
$ws.Range("J85").Value = 24997.143
$ws.Range("L85").Value = 24997.143
$ws.Range("N85").Value = -27649.143
$ws.Range("H99").Value = 3073.3333
$ws.Range("I99").Value = 1980
$ws.Range("K99").Value = 1980
$ws.Range("M99").Value = -482

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1152.5
$ws.Range("I16").Value = 904.5
$ws.Range("K16").Value = 904.5
$ws.Range("M16").Value = -617.5
$ws.Range("H31").Value = 1596.3529
$ws.Range("I31").Value = 1630.4828
$ws.Range("K31").Value = 1630.4828
$ws.Range("M31").Value = -1335.4828
$ws.Range("H34").Value = 1596.3529
$ws.Range("I34").Value = 1630.4828
$ws.Range("K34").Value = 1630.4828
$ws.Range("M34").Value = -1428.4828
$ws.Range("H92").Value = 28639.8
$ws.Range("J92").Value = 28639.8
$ws.Range("L92").Value = 28639.8
$ws.Range("N92").Value = -33631.8
$ws.Range("H113").Value = 1152.5
$ws.Range("I113").Value = 904.5
$ws.Range("K113").Value = 904.5
$ws.Range("M113").Value = 1265.5
$ws.Range("H132").Value = 2472.913
$ws.Range("I132").Value = 2279.1904
$ws.Range("K132").Value = 6837.5712
$ws.Range("M132").Value = -4307.5712
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 102447.35
$ws.Range("I129").Value = 250333.12
$ws.Range("J129").Value = 3856.8333
$ws.Range("K129").Value = 750999.36
$ws.Range("L129").Value = 11570.4999
$ws.Range("M129").Value = -745999.36
$ws.Range("N129").Value = -21570.4999
$ws.Range("H131").Value = 102846.5
$ws.Range("J131").Value = 3187.125
$ws.Range("L131").Value = 9561.375
$ws.Range("N131").Value = -19641.375

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 17890.908
$ws.Range("I99").Value = 13845.375
$ws.Range("K99").Value = 13845.375
$ws.Range("M99").Value = -11599.375
$ws.Range("H107").Value = 33335556
$ws.Range("I107").Value = 868.875
$ws.Range("J107").Value = 71432344
$ws.Range("K107").Value = 868.875
$ws.Range("L107").Value = 71432344
$ws.Range("M107").Value = 1051.125
$ws.Range("N107").Value = -71436184
$ws.Range("H113").Value = 25002256
$ws.Range("I113").Value = 29413770
$ws.Range("J113").Value = 3669.3333
$ws.Range("K113").Value = 29413770
$ws.Range("L113").Value = 3669.3333
$ws.Range("M113").Value = -29411600
$ws.Range("N113").Value = -8009.3333
$ws.Range("H135").Value = 172889.5
$ws.Range("J135").Value = 172889.5
$ws.Range("L135").Value = 172889.5
$ws.Range("N135").Value = -183029.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1141.5555
$ws.Range("J22").Value = 1161.75
$ws.Range("L22").Value = 1161.75
$ws.Range("N22").Value = -1751.75
$ws.Range("H27").Value = 1141.5555
$ws.Range("J27").Value = 1161.75
$ws.Range("L27").Value = 1161.75
$ws.Range("N27").Value = -1375.75
$ws.Range("H61").Value = 633
$ws.Range("J61").Value = 699.5
$ws.Range("L61").Value = 699.5
$ws.Range("N61").Value = -1103.5
$ws.Range("H98").Value = 14451.667
$ws.Range("J98").Value = 14451.667
$ws.Range("L98").Value = 14451.667
$ws.Range("N98").Value = -20441.667
$ws.Range("H113").Value = 633
$ws.Range("J113").Value = 699.5
$ws.Range("L113").Value = 699.5
$ws.Range("N113").Value = -5039.5
$ws.Range("H136").Value = 4299.091
$ws.Range("I136").Value = 3359
$ws.Range("K136").Value = 10077
$ws.Range("M136").Value = -7527

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 3005000
$ws.Range("J14").Value = 10000
$ws.Range("L14").Value = 10000
$ws.Range("N14").Value = -10336
$ws.Range("H107").Value = 100001140
$ws.Range("I107").Value = 1630.6666
$ws.Range("J107").Value = 250000400
$ws.Range("K107").Value = 4891.9998
$ws.Range("L107").Value = 750001200
$ws.Range("M107").Value = -2971.9998
$ws.Range("N107").Value = -750005040
$ws.Range("H122").Value = 2011.2222
$ws.Range("I122").Value = 1482.4
$ws.Range("K122").Value = 4447.200000000001
$ws.Range("M122").Value = -1997.200000000001
$ws.Range("H136").Value = 884.4211
$ws.Range("I136").Value = 884.4211
$ws.Range("K136").Value = 2653.2633
$ws.Range("M136").Value = -103.2633000000001
